# Update imputed values in the RandomForest result sheet (Sheet1).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A3").Value = -22.00050000000001
$ws.Range("A21").Value = -20.08209999999999
$ws.Range("A23").Value = -20.74159999999998
$ws.Range("A25").Value = -21.53059999999999
$ws.Range("E27").Value = 16.75799999999999
$ws.Range("E31").Value = 16.73309999999999
$ws.Range("E39").Value = 15.9079
$ws.Range("E48").Value = 17.4302
$ws.Range("E51").Value = 17.385
$ws.Range("E52").Value = 17.2761
$ws.Range("A53").Value = -21.89029999999999
$ws.Range("E55").Value = 16.5608
$ws.Range("E56").Value = 16.4408
$ws.Range("A57").Value = -22.65220000000001
$ws.Range("E57").Value = 16.6427
$ws.Range("A59").Value = -22.3629
$ws.Range("A69").Value = -21.62659999999998
$ws.Range("E73").Value = 17.2654
$ws.Range("A79").Value = -20.46520000000002
$ws.Range("A83").Value = -21.822
$ws.Range("E89").Value = 17.23140000000001
$ws.Range("E90").Value = 16.56109999999999
$ws.Range("A93").Value = -21.2698
